# Insert a new price record as row 228 (Ñuble / Terminal Hortofrutícola Agro
# Chillán - Zanahoria, weekly price update), pushing existing rows 228-318
# down to 229-319.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(228).Insert()

$ws.Range("A228").Value = 7
$ws.Range("B228").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C228").Value = 'Ñuble'
$ws.Range("D228").Value2 = 44795
$ws.Range("E228").Value = 16
$ws.Range("F228").Value = 100114013
$ws.Range("G228").Value = 'Zanahoria'
$ws.Range("H228").Value = 'Sin especificar'
$ws.Range("I228").Value = 'Primera'
$ws.Range("J228").Value = 120
$ws.Range("K228").Value = 10000
$ws.Range("L228").Value = 11000
$ws.Range("M228").Value = 10500
$ws.Range("N228").Value = '$/saco 20 kilos'
$ws.Range("O228").Value = 'Provincia de Diguillín'
$ws.Range("P228").Value = 525
$ws.Range("Q228").Value = 20
$ws.Range("R228").Value = 'Hortaliza'
